$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.555.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.835.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.87%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4241'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3681'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07255'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8670'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.26%  '

$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.017.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.37%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.378'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.505'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06970'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.09%  '

$ws.Range("E16").Value = '  +1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.65'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009009'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.010'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.039.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.033'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.71%  '

$ws.Range("E24").Value = '  +3.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.974'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.247'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.829'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08898'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7685'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.539'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.960'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.142'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.013'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.25%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.103'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05363'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01941'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.827'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5094'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1659'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.58%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.783'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.466'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.83%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.12%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06541'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '105.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.010'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4677'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.618'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.795'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.82%  '

